# Update cryptos list with latest scraped price / volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    # Force the cell to hold the literal text even when it looks like a
    # plain number (Excel would otherwise silently convert it to a
    # numeric value, e.g. "601.30" -> 601.3). Resetting the style back to
    # Normal afterwards keeps the cell's formatting identical to before.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.772.03"
$ws.Range("E2").Value = "  +1.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.742.95"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "601.30"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6 - Solana
Set-TextValue "D6" "168.96"
$ws.Range("E6").Value = "  -1.68%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.742.78"
$ws.Range("E7").Value = "  -1.52%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.81%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.87%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.44%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.462"
$ws.Range("E12").Value = "  -1.21%  "

# Row 13 - Avalanche
Set-TextValue "D13" "38.19"
$ws.Range("E13").Value = "  -1.58%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +0.50%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.363.08"
$ws.Range("E15").Value = "  -1.79%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.736.81"
$ws.Range("E16").Value = "  -1.80%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "68.756.93"
$ws.Range("E17").Value = "  +1.26%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +0.57%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -0.30%  "

# Row 20 - Chainlink
Set-TextValue "D20" "17.13"
$ws.Range("E20").Value = "  -0.88%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.78"
$ws.Range("E21").Value = "  +16.49%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "494.81"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -1.80%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.39"
$ws.Range("E24").Value = "  -0.43%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +0.26%  "

# Row 26 - Fetch.AI
Set-TextValue "D26" "2.31"
$ws.Range("E26").Value = "  -3.64%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +0.30%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -0.62%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.09%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +4.27%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.11%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  +0.70%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "31.89"
$ws.Range("E33").Value = "  -3.93%  "

# Row 34 - WrappedeETH
$ws.Range("D34").Value = "3.884.11"
$ws.Range("E34").Value = "  -1.53%  "

# Row 35 - RenzoRestakedETH
$ws.Range("D35").Value = "3.671.91"
$ws.Range("E35").Value = "  -1.79%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  -1.60%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - Mantle
$ws.Range("E38").Value = "  -0.80%  "

# Row 39 - Filecoin
Set-TextValue "D39" "5.84"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -0.03%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -1.72%  "

# Row 42 - Bittensor
Set-TextValue "D42" "437.52"
$ws.Range("E42").Value = "  -5.60%  "

# Row 43 - OKB
Set-TextValue "D43" "48.92"
$ws.Range("E43").Value = "  -0.34%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -1.62%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  +0.22%  "

# Row 46 - Cosmos
$ws.Range("E46").Value = "  +0.84%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  +0.01%  "

# Row 48 - Arweave
$ws.Range("E48").Value = "  +0.13%  "

# Row 49 - was Monero, now Maker (rows 49 and 50 swapped their coin identity)
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.809.24"
$ws.Range("E49").Value = "  -1.25%  "

# Row 50 - was Maker, now Monero
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "141.15"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +0.74%  "
